# Estimated CPI from years from 1920 - 1970 for which CPI data are not available
# Inflation rate is assumed to be average of the next 5 years; and CPI is
# calculated based on inflation rate using standard inflation fomular.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")
$wsConst = $wb.Worksheets.Item("country_constants")

# 1. Insert 4 new columns before column M (1971-1974), shifting 1975.. to the right.
$ws.Columns("M:P").Insert()

# 2. Fill in the new year header values in row 1.
$ws.Range("M1").Value = 1971
$ws.Range("N1").Value = 1972
$ws.Range("O1").Value = 1973
$ws.Range("P1").Value = 1974

# 3. Row 37 - inflation rate. Update the early (pre-1970) placeholder values
#    and populate the newly inserted columns with estimated inflation rates.
$ws.Range("E37:K37").Value = 0.123

$ws.Range("M37").Value = 0.091
$ws.Range("N37").Value = 0.22
$ws.Range("O37").Value = 0.111
$ws.Range("P37").Value = 0.145

# Give the newly inserted inflation-rate cells the same look as the rest of
# the row (plain font, vertically centered) by copying format from a cell
# that already uses that font, then centering it vertically.
$wsConst.Range("C2").Copy()
$ws.Range("M37:P37").PasteSpecial(-4122)
$ws.Range("M37:P37").VerticalAlignment = -4108

# 4. Row 38 - CPI values. Update the pre-1970 placeholder values and the
#    values shifted from the old K/L columns, then populate the newly
#    inserted columns with the calculated CPI.
$ws.Range("E38").Value = 4.31
$ws.Range("F38").Value = 4.84
$ws.Range("G38").Value = 5.44
$ws.Range("H38").Value = 6.11
$ws.Range("I38").Value = 6.86
$ws.Range("J38").Value = 7.7
$ws.Range("K38").Value = 8.68
$ws.Range("L38").Value = 9.04

$ws.Range("M38").Value = 9.41
$ws.Range("N38").Value = 11.48
$ws.Range("O38").Value = 12.75
$ws.Range("P38").Value = 14.6

# Newly inserted CPI cells should use the same number style as the
# neighbouring K38/L38/Q38 cells.
$ws.Range("L38").Copy()
$ws.Range("M38:P38").PasteSpecial(-4122)

# 5. Restore the application's active cell/selection on the sheet.
$ws.Range("H45").Select()

Write-Host "done"
